$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Notification System"
$ws.Range("B2").Value = "Done"
$ws.Range("A3").Value = "URL Shortner"
$ws.Range("B3").Value = "Done"

$ws.Columns.Item(1).ColumnWidth = 40.5

$ws.Range("B3").Select()
